$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing "Kosten in €" header (currently in G2) up into a new G1 header cell,
# and set the real price values (as text, with comma decimal separators) in G2:G8.
$ws.Range("G1").Value = "Kosten in €"

$ws.Range("G2").Value = "10,96"
$ws.Range("G3").Value = "5,48"
$ws.Range("G4").Value = "43,84"
$ws.Range("G5").Value = "67,20"
$ws.Range("G6").Value = "16,80"
$ws.Range("G7").Value = "67,20"
$ws.Range("G8").Value = "33,60"
